# Rename "LATAM" sheets to "LAC" (Latin America and the Caribbean), matching
# the refreshed supplementary-information workbook for the new journal article.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("imp_latam").Name = "imp_lac"
$wb.Worksheets.Item("imp_latam_all").Name = "imp_lac_all"
$wb.Worksheets.Item("imp_latam_sum").Name = "imp_lac_sum"

# The pie/bar charts on those sheets keep their own series formulas, which
# are not retargeted automatically by a sheet rename - fix them up so they
# keep pointing at the (renamed) source ranges.
foreach ($ws in $wb.Worksheets) {
    foreach ($co in $ws.ChartObjects()) {
        $chart = $co.Chart
        foreach ($ser in $chart.SeriesCollection()) {
            $f = $ser.Formula
            $newf = $f -replace "imp_latam_sum", "imp_lac_sum"
            $newf = $newf -replace "imp_latam", "imp_lac"
            if ($newf -ne $f) {
                $ser.Formula = $newf
            }
        }
    }
}

# Refresh the "Cover" sheet text for the new supplementary-information
# release (new title, new reference citation, new update date, renamed
# guide entries, and the LATAM -> LAC abbreviation).
$ws = $wb.Worksheets.Item("Cover")
$ws.Range("A1").Value = "Supplementary Information"
$ws.Range("A4").Value = "This supplementary information includes the impact analysis of trade-offs and synergies from the circular economy scenarios proposed by 'How to measure Circularity Trade-offs and Synergies?'"
$ws.Range("C7").Value = "November 25, 2024"
$ws.Range("B14").Value = "imp_lac"
$ws.Range("B15").Value = "imp_lac_all"
$ws.Range("B16").Value = "imp_lac_sum"
$ws.Range("D21").Value = "LAC"

# Match the saved selection on the Cover sheet.
$ws.Activate()
$ws.Range("A5").Select()
